# Apply the "new LM training" metrics update to Sheet1.
# - Column A (model name) gets reordered/relabeled per row 2..26
# - Columns B..Q get a single, constant set of new metric values for every row 2..26

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New model-name labels for rows 2..26 (A2..A26)
$names = @(
    "model_20_5_0",
    "model_20_5_22",
    "model_20_5_21",
    "model_20_5_20",
    "model_20_5_19",
    "model_20_5_18",
    "model_20_5_17",
    "model_20_5_16",
    "model_20_5_15",
    "model_20_5_14",
    "model_20_5_13",
    "model_20_5_23",
    "model_20_5_12",
    "model_20_5_10",
    "model_20_5_9",
    "model_20_5_8",
    "model_20_5_7",
    "model_20_5_6",
    "model_20_5_5",
    "model_20_5_4",
    "model_20_5_3",
    "model_20_5_2",
    "model_20_5_1",
    "model_20_5_11",
    "model_20_5_24"
)

# New constant metric values for columns B..Q (r2, r2_sup, r2_test, r2_val, r2_vt,
# mse, mse_sup, mse_test, mse_val, mse_vt, mape, rmse, r2_adj, rsd, aic, bic)
# Parsed via [double] cast from strings so scientific notation (e.g. 1e-05) is
# handled reliably regardless of numeric-literal parsing quirks.
$values = @(
    [double]"0.9999805300082542",
    [double]"0.9991182316315311",
    [double]"0.9999999999999697",
    [double]"0.9999968636398666",
    [double]"0.9999998123027997",
    [double]"1.817438266093276e-05",
    [double]"0.0008230920668076507",
    [double]"3.311478096098398e-14",
    [double]"3.066770497431888e-07",
    [double]"1.533385414289849e-07",
    [double]"0.0002731515751168351",
    [double]"0.004263142345844525",
    [double]"1.000035944600146",
    [double]"0.004444633141428372",
    [double]"95.83099500301532",
    [double]"140.9294005231387"
)

for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $i + 2

    # Column A: update model name label
    $ws.Cells.Item($row, 1).Value = $names[$i]

    # Columns B..Q (2..17): update metric values
    for ($j = 0; $j -lt $values.Count; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
